$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'286.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'2.40%"
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'28.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'4.09%"
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("D4").Value = "'5.069"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'4.84%"
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'0.06643"
$ws.Range("D5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'7.360"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'4.59%"
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "'3.399"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'2.14%"
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").Value = "'1.372"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'4.26%"
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9385"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'3.86%"
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1567"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'1.91%"
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.06617"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'7.82%"
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.07601"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'1.25%"
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.02944"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.50%"
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.08981"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.08%"
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001591"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.88%"
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D16").Value = "'0.04504"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'2.12%"
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("B17").Value = "One"
$ws.Range("C17").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D17").Value = "'0.0006449"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.36%"
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").Value = "'0.006284"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'3.62%"
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("B19").Value = "LEO"
$ws.Range("C19").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D19").Value = "'3.440"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-1.49%"
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'2.251"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'1.16%"
$ws.Range("E20").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'0.1297"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-4.03%"
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'4.067"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'3.86%"
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'0.1553"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'3.27%"
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'0.001183"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.68%"
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("E26").Value = "'4.26%"
$ws.Range("E26").Style = "Normal"

# Row 28
$ws.Range("D28").Value = "'0.0001618"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'-2.31%"
$ws.Range("E28").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'0.04196"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'3.15%"
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("D41").Value = "'0.006751"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'1.40%"
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("D42").Value = "'0.1248"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-10.43%"
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = "'0.002021"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-3.28%"
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'0.01232"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'11.62%"
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'0.00005611"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'1.17%"
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("E46").Value = "'20.74%"
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'0.01307"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-29.27%"
$ws.Range("E47").Style = "Normal"
